# "wrapping up test file audit"
#
# The workbook had a stray leftover row (A16="Sheet", B16=3, C16=4) on the
# "optimization_parameters" sheet -- apparently debris from an earlier
# edit/paste. Remove it (which shifts the "simulation_timepoints" row up
# from row 17 to row 16), then leave the review focus on the
# "optimization_diagnostics" sheet (the last sheet that was checked during
# the audit).

$wb = $excel.ActiveWorkbook

# --- Clean up the stray "Sheet" row on optimization_parameters ---
$wsParams = $wb.Worksheets.Item("optimization_parameters")
$wsParams.Activate()
$wsParams.Rows.Item(16).Delete()
# Leave the selection on the row that slid up into the gap, matching what
# Excel leaves selected after a row delete.
$wsParams.Rows.Item(16).Select()

# --- Move the active tab to optimization_diagnostics (last sheet reviewed) ---
$wsDiag = $wb.Worksheets.Item("optimization_diagnostics")
$wsDiag.Activate()
